$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 113, shifting existing rows 113:128 down to 114:129
$ws.Rows.Item(113).Insert()

$ws.Range("A113").Value = 5
$ws.Range("B113").Value = "Macroferia Regional de Talca"
$ws.Range("C113").Value = "Maule"
$ws.Range("D113").Value = 45194
$ws.Range("E113").Value = 7
$ws.Range("F113").Value = 100112040
$ws.Range("G113").Value = "Cilantro"
$ws.Range("H113").Value = "Sin especificar"
$ws.Range("I113").Value = "Primera"
$ws.Range("J113").Value = 300
$ws.Range("K113").Value = 8000
$ws.Range("L113").Value = 8000
$ws.Range("M113").Value = 8000
$ws.Range("N113").Value = "$/caja 36 atados"
$ws.Range("O113").Value = "Región Metropolitana"
$ws.Range("P113").Value = 222
$ws.Range("Q113").Value = 36
$ws.Range("R113").Value = "Hortaliza"
